$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1766.5834
$ws.Range("J2").Value = 1216.6666
$ws.Range("L2").Value = 1216.6666
$ws.Range("N2").Value = -1442.6666
$ws.Range("H51").Value = 5091.5
$ws.Range("I51").Value = 7400
$ws.Range("K51").Value = 7400
$ws.Range("M51").Value = -6916
$ws.Range("H57").Value = 94890
$ws.Range("J57").Value = 94890
$ws.Range("L57").Value = 284670
$ws.Range("N57").Value = -285668
$ws.Range("H62").Value = 4340.2
$ws.Range("I62").Value = 3733.5557
$ws.Range("K62").Value = 3733.5557
$ws.Range("M62").Value = -3109.5557
$ws.Range("H64").Value = 7744.3335
$ws.Range("J64").Value = 9000
$ws.Range("L64").Value = 9000
$ws.Range("N64").Value = -9496
$ws.Range("H65").Value = 4340.2
$ws.Range("I65").Value = 3733.5557
$ws.Range("K65").Value = 18667.7785
$ws.Range("M65").Value = -15547.7785
$ws.Range("H67").Value = 7744.3335
$ws.Range("J67").Value = 9000
$ws.Range("L67").Value = 9000
$ws.Range("N67").Value = -10716
$ws.Range("H125").Value = 2342.9148
$ws.Range("I125").Value = 1464.4
$ws.Range("K125").Value = 13179.6
$ws.Range("M125").Value = -10719.6

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 24999.5
$ws.Range("I33").Value = 24999
$ws.Range("J33").Value = 25000
$ws.Range("K33").Value = 24999
$ws.Range("L33").Value = 25000
$ws.Range("M33").Value = -24670
$ws.Range("N33").Value = -25658
$ws.Range("H37").Value = 44745
$ws.Range("I37").Value = 44745
$ws.Range("K37").Value = 44745
$ws.Range("M37").Value = -44472
$ws.Range("H55").Value = 49999
$ws.Range("J55").Value = 49999
$ws.Range("L55").Value = 49999
$ws.Range("N55").Value = -50629
$ws.Range("H63").Value = 2877.4
$ws.Range("I63").Value = 2621.75
$ws.Range("K63").Value = 2621.75
$ws.Range("M63").Value = -1935.75
$ws.Range("H66").Value = 2877.4
$ws.Range("I66").Value = 2621.75
$ws.Range("K66").Value = 13108.75
$ws.Range("M66").Value = -9676.75
$ws.Range("H88").Value = 2178.75
$ws.Range("J88").Value = 3676.25
$ws.Range("L88").Value = 3676.25
$ws.Range("N88").Value = -4488.25
$ws.Range("H91").Value = 2178.75
$ws.Range("J91").Value = 3676.25
$ws.Range("L91").Value = 3676.25
$ws.Range("N91").Value = -6484.25
$ws.Range("H122").Value = 17097908
$ws.Range("I122").Value = 27780788
$ws.Range("K122").Value = 83342364
$ws.Range("M122").Value = -83339914

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2159.4
$ws.Range("I107").Value = 1669.8182
$ws.Range("J107").Value = 3505.75
$ws.Range("K107").Value = 1669.8182
$ws.Range("L107").Value = 3505.75
$ws.Range("M107").Value = 250.1818000000001
$ws.Range("N107").Value = -7345.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2999.75
$ws.Range("J62").Value = 2999.75
$ws.Range("L62").Value = 2999.75
$ws.Range("N62").Value = -4247.75
$ws.Range("H65").Value = 2999.75
$ws.Range("J65").Value = 2999.75
$ws.Range("L65").Value = 14998.75
$ws.Range("N65").Value = -21238.75
$ws.Range("H139").Value = 76654.5
$ws.Range("I139").Value = 13309
$ws.Range("K139").Value = 13309
$ws.Range("M139").Value = -8169

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1096.6154
$ws.Range("I5").Value = 710.2353000000001
$ws.Range("J5").Value = 1826.4445
$ws.Range("K5").Value = 2130.7059
$ws.Range("L5").Value = 5479.333500000001
$ws.Range("M5").Value = -2018.7059
$ws.Range("N5").Value = -5703.333500000001
$ws.Range("H86").Value = 200
$ws.Range("I86").Value = 200
$ws.Range("J86").Value = 200
$ws.Range("K86").Value = 600
$ws.Range("L86").Value = 600
$ws.Range("M86").Value = 586
$ws.Range("N86").Value = -2972
$ws.Range("H89").Value = 200
$ws.Range("I89").Value = 200
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 1800
$ws.Range("L89").Value = 1800
$ws.Range("M89").Value = 4128
$ws.Range("N89").Value = -13656
$ws.Range("H98").Value = 785.7778
$ws.Range("I98").Value = 754.8
$ws.Range("J98").Value = 824.5
$ws.Range("K98").Value = 2264.4
$ws.Range("L98").Value = 2473.5
$ws.Range("M98").Value = -766.3999999999996
$ws.Range("N98").Value = -5469.5
$ws.Range("H105").Value = 4828.6
$ws.Range("J105").Value = 4828.6
$ws.Range("L105").Value = 14485.8
$ws.Range("N105").Value = -19727.8
$ws.Range("H132").Value = 1722.8064
$ws.Range("J132").Value = 1970.5454
$ws.Range("L132").Value = 17734.9086
$ws.Range("N132").Value = -22794.9086
$ws.Range("H135").Value = 1096.6154
$ws.Range("I135").Value = 710.2353000000001
$ws.Range("J135").Value = 1826.4445
$ws.Range("K135").Value = 6392.117700000001
$ws.Range("L135").Value = 16438.0005
$ws.Range("M135").Value = -3857.117700000001
$ws.Range("N135").Value = -21508.0005

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1050.6666
$ws.Range("I107").Value = 1050.6666
$ws.Range("K107").Value = 1050.6666
$ws.Range("M107").Value = 869.3334
$ws.Range("H122").Value = 251220.92
$ws.Range("I122").Value = 389624.97
$ws.Range("J122").Value = 6352.231
$ws.Range("K122").Value = 1168874.91
$ws.Range("L122").Value = 19056.693
$ws.Range("M122").Value = -1166424.91
$ws.Range("N122").Value = -23956.693

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4290.6665
$ws.Range("I46").Value = 3173.1538
$ws.Range("J46").Value = 6106.625
$ws.Range("K46").Value = 3173.1538
$ws.Range("L46").Value = 6106.625
$ws.Range("M46").Value = -2985.1538
$ws.Range("N46").Value = -6482.625
$ws.Range("H132").Value = 6299.8
$ws.Range("I132").Value = 6693.067
$ws.Range("K132").Value = 20079.201
$ws.Range("M132").Value = -17549.201
$ws.Range("H136").Value = 60478.39
$ws.Range("I136").Value = 87051
$ws.Range("J136").Value = 7333.1665
$ws.Range("K136").Value = 261153
$ws.Range("L136").Value = 21999.4995
$ws.Range("M136").Value = -258603
$ws.Range("N136").Value = -27099.4995

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1464.2858
$ws.Range("I81").Value = 1458.3334
$ws.Range("K81").Value = 2916.6668
$ws.Range("M81").Value = -1855.6668
$ws.Range("H84").Value = 1464.2858
$ws.Range("I84").Value = 1458.3334
$ws.Range("K84").Value = 14583.334
$ws.Range("M84").Value = -9279.333999999999
$ws.Range("H136").Value = 2870.7334
$ws.Range("I136").Value = 2222.7827
$ws.Range("J136").Value = 4999.7144
$ws.Range("K136").Value = 6668.348100000001
$ws.Range("L136").Value = 14999.1432
$ws.Range("M136").Value = -4118.348100000001
$ws.Range("N136").Value = -20099.1432
